$wb = $excel.ActiveWorkbook

# StatOutput sheet: numbers went from files=1/sample=2 to files=0/sample=0
# (stored as text, matching the original file's shared-string typed cells)
$statOutput = $wb.Worksheets.Item("StatOutput")
$statOutput.Range("A2").NumberFormat = "@"
$statOutput.Range("B2").NumberFormat = "@"
$statOutput.Range("A2").Value = "0"
$statOutput.Range("B2").Value = "0"

# StatOutput_Message sheet: the logged cypher query (row 18, col A) is updated
# to reflect the breed filter actually used (Cavalier King Charles Spaniel)
# together with the OPTIONAL MATCH / count(...) style query.
$statMessage = $wb.Worksheets.Item("StatOutput_Message")
$newQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Cavalier King Charles Spaniel']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$statMessage.Range("A18").Value = $newQuery
